# Schémas.xlsx — mise à jour schéma + correction oubli + nettoyage
#
# - Mise à jour du schéma avec l'augmentation de proba pour les ressources doublées
# - Correction de l'oubli de cette modification pour les plaines
# - Retrait des commentaires signalant une modification (cellule "s=27" déplacée
#   en fin de table de styles lors du nettoyage — effet de bord cosmétique sans
#   incidence sur la mise en forme visible)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ressources")

# --- Mise à jour des probabilités (ressources doublées) ---
$ws.Range("D4").Value = 0.2

$ws.Range("C5").Value = 0.3
$ws.Range("D5").Value = 0.4
$ws.Range("E5").Value = 0.5
$ws.Range("K5").Value = 0.25

$ws.Range("C6").Value = 0.35

$ws.Range("E8").Value = 0.25
$ws.Range("K8").Value = 0.25

# --- Correction de l'oubli pour les plaines ---
$ws.Range("D29").Value = 0.35

# --- Nettoyage de la vue (on se replace sur la cellule D11) ---
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
